$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 147
$ws.Range("F5").Value = 59
$ws.Range("F6").Value = 3728
$ws.Range("F8").Value = 2506
$ws.Range("F10").Value = 2993
$ws.Range("F12").Value = 528
$ws.Range("F13").Value = 2271
$ws.Range("G13").Value = 65
$ws.Range("F15").Value = 111
$ws.Range("F16").Value = 73
$ws.Range("F17").Value = 428
$ws.Range("F20").Value = 333
$ws.Range("F21").Value = 293
$ws.Range("F22").Value = 316
$ws.Range("F23").Value = 632
$ws.Range("F24").Value = 1374
$ws.Range("F26").Value = 1285
$ws.Range("F27").Value = 120
$ws.Range("F28").Value = 143
$ws.Range("F29").Value = 241
$ws.Range("F30").Value = 4119
$ws.Range("F31").Value = 3710
$ws.Range("F32").Value = 60
$ws.Range("F34").Value = 1094
$ws.Range("F35").Value = 445
$ws.Range("F37").Value = 1301
$ws.Range("F38").Value = 141
$ws.Range("F43").Value = 45

# Sheet: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 186
$ws.Range("F8").Value = 17
$ws.Range("F16").Value = 196

# Sheet: 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 1018
$ws.Range("F4").Value = 2210

# Sheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 1018
$ws.Range("F7").Value = 147
$ws.Range("F8").Value = 186
$ws.Range("F9").Value = 59
$ws.Range("F11").Value = 3728
$ws.Range("F13").Value = 2506
$ws.Range("F15").Value = 2993
$ws.Range("F16").Value = 528
$ws.Range("F17").Value = 2271
$ws.Range("G17").Value = 65
$ws.Range("F19").Value = 111
$ws.Range("F20").Value = 73
$ws.Range("F21").Value = 428
$ws.Range("F23").Value = 333
$ws.Range("F24").Value = 316
$ws.Range("F25").Value = 632
$ws.Range("F26").Value = 1374
$ws.Range("F28").Value = 1285
$ws.Range("F29").Value = 120
$ws.Range("F30").Value = 143
$ws.Range("F32").Value = 17
$ws.Range("F33").Value = 4119
$ws.Range("F34").Value = 3710
$ws.Range("F35").Value = 60
$ws.Range("F38").Value = 445
$ws.Range("F43").Value = 1301
$ws.Range("F44").Value = 141
$ws.Range("F48").Value = 45
$ws.Range("F49").Value = 196
